# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Gilgamesh_Profits workbook sheets
# as described by the source diff (scheduled runner market-data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 866.3333
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -1225.5
$ws.Range("H8").Value = 1732.6666
$ws.Range("I8").Value = 99
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 297
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = -158
$ws.Range("N8").Value = -15278
$ws.Range("H18").Value = 125001224
$ws.Range("I18").Value = 892.75
$ws.Range("J18").Value = 250001550
$ws.Range("K18").Value = 892.75
$ws.Range("L18").Value = 250001550
$ws.Range("M18").Value = -608.75
$ws.Range("N18").Value = -250002118
$ws.Range("H43").Value = 1800.4736
$ws.Range("I43").Value = 1727.1
$ws.Range("K43").Value = 1727.1
$ws.Range("M43").Value = -1658.1
$ws.Range("H54").Value = 10027.667
$ws.Range("I54").Value = 4999.5
$ws.Range("J54").Value = 20084
$ws.Range("K54").Value = 4999.5
$ws.Range("L54").Value = 20084
$ws.Range("M54").Value = -4513.5
$ws.Range("N54").Value = -21056
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 17583.334
$ws.Range("J74").Value = 10000
$ws.Range("L74").Value = 10000
$ws.Range("N74").Value = -11872
$ws.Range("H77").Value = 17583.334
$ws.Range("J77").Value = 10000
$ws.Range("L77").Value = 50000
$ws.Range("N77").Value = -59360
$ws.Range("H94").Value = 1072.25
$ws.Range("I94").Value = 1030.3334
$ws.Range("J94").Value = 1198
$ws.Range("K94").Value = 1030.3334
$ws.Range("L94").Value = 1198
$ws.Range("M94").Value = -579.3334
$ws.Range("N94").Value = -2100
$ws.Range("H141").Value = 2143.7693
$ws.Range("I141").Value = 2122.15
$ws.Range("J141").Value = 2215.8333
$ws.Range("K141").Value = 6366.450000000001
$ws.Range("L141").Value = 6647.499899999999
$ws.Range("M141").Value = -1186.450000000001
$ws.Range("N141").Value = -17007.4999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2106
$ws.Range("I63").Value = 2106
$ws.Range("K63").Value = 2106
$ws.Range("M63").Value = -1420
$ws.Range("H66").Value = 2106
$ws.Range("I66").Value = 2106
$ws.Range("K66").Value = 10530
$ws.Range("M66").Value = -7098
$ws.Range("H74").Value = 2717.72
$ws.Range("I74").Value = 2467.5293
$ws.Range("J74").Value = 3249.375
$ws.Range("K74").Value = 2467.5293
$ws.Range("L74").Value = 3249.375
$ws.Range("M74").Value = -1593.5293
$ws.Range("N74").Value = -4997.375
$ws.Range("H76").Value = 72000
$ws.Range("J76").Value = 72000
$ws.Range("L76").Value = 72000
$ws.Range("N76").Value = -72676
$ws.Range("H77").Value = 2717.72
$ws.Range("I77").Value = 2467.5293
$ws.Range("J77").Value = 3249.375
$ws.Range("K77").Value = 12337.6465
$ws.Range("L77").Value = 16246.875
$ws.Range("M77").Value = -7969.646500000001
$ws.Range("N77").Value = -24982.875
$ws.Range("H79").Value = 72000
$ws.Range("J79").Value = 72000
$ws.Range("L79").Value = 72000
$ws.Range("N79").Value = -74340
$ws.Range("H122").Value = 2338.8096
$ws.Range("I122").Value = 2338.8096
$ws.Range("K122").Value = 7016.4288
$ws.Range("M122").Value = -4566.4288
$ws.Range("H132").Value = 3009.16
$ws.Range("I132").Value = 2826.875
$ws.Range("K132").Value = 8480.625
$ws.Range("M132").Value = -5950.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17247808
$ws.Range("I20").Value = 20840436
$ws.Range("K20").Value = 20840436
$ws.Range("M20").Value = -20840189

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 16669076
$ws.Range("I62").Value = 20002294
$ws.Range("J62").Value = 2992
$ws.Range("K62").Value = 20002294
$ws.Range("L62").Value = 2992
$ws.Range("M62").Value = -20001670
$ws.Range("N62").Value = -4240
$ws.Range("H65").Value = 16669076
$ws.Range("I65").Value = 20002294
$ws.Range("J65").Value = 2992
$ws.Range("K65").Value = 100011470
$ws.Range("L65").Value = 14960
$ws.Range("M65").Value = -100008350
$ws.Range("N65").Value = -21200
$ws.Range("H105").Value = 1804.9166
$ws.Range("I105").Value = 1207.375
$ws.Range("K105").Value = 1207.375
$ws.Range("M105").Value = 539.625
$ws.Range("H134").Value = 2704.0833
$ws.Range("I134").Value = 2481.6667
$ws.Range("J134").Value = 3816.1667
$ws.Range("K134").Value = 7445.000100000001
$ws.Range("L134").Value = 11448.5001
$ws.Range("M134").Value = -4910.000100000001
$ws.Range("N134").Value = -16518.5001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 371.25
$ws.Range("I33").Value = 359.5
$ws.Range("J33").Value = 383
$ws.Range("K33").Value = 2157
$ws.Range("L33").Value = 2298
$ws.Range("M33").Value = -1874
$ws.Range("N33").Value = -2864
$ws.Range("H80").Value = 15800
$ws.Range("J80").Value = 19850
$ws.Range("L80").Value = 59550
$ws.Range("N80").Value = -61422
$ws.Range("H83").Value = 15800
$ws.Range("J83").Value = 19850
$ws.Range("L83").Value = 178650
$ws.Range("N83").Value = -188010
$ws.Range("H93").Value = 7332.8335
$ws.Range("I93").Value = 5999
$ws.Range("K93").Value = 17997
$ws.Range("M93").Value = -16125
$ws.Range("H103").Value = 2663
$ws.Range("I103").Value = 326.5
$ws.Range("K103").Value = 979.5
$ws.Range("M103").Value = -100.5
$ws.Range("H131").Value = 5441801.5
$ws.Range("I131").Value = 16684128
$ws.Range("J131").Value = 1966.3226
$ws.Range("K131").Value = 50052384
$ws.Range("L131").Value = 5898.9678
$ws.Range("M131").Value = -50047344
$ws.Range("N131").Value = -15978.9678

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 999.7857
$ws.Range("I22").Value = 521.6
$ws.Range("J22").Value = 1265.4445
$ws.Range("K22").Value = 521.6
$ws.Range("L22").Value = 1265.4445
$ws.Range("M22").Value = 7.399999999999977
$ws.Range("N22").Value = -2323.4445
$ws.Range("H36").Value = 4017
$ws.Range("I36").Value = 4017
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4017
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = -3532
$ws.Range("N36").Value = 0
$ws.Range("H70").Value = 7198.2856
$ws.Range("I70").Value = 4996.3335
$ws.Range("K70").Value = 4996.3335
$ws.Range("M70").Value = -4726.3335
$ws.Range("H73").Value = 7198.2856
$ws.Range("I73").Value = 4996.3335
$ws.Range("K73").Value = 4996.3335
$ws.Range("M73").Value = -4060.3335
$ws.Range("H80").Value = 90913570
$ws.Range("I80").Value = 166670370
$ws.Range("J80").Value = 5418.8
$ws.Range("K80").Value = 166670370
$ws.Range("L80").Value = 5418.8
$ws.Range("M80").Value = -166669372
$ws.Range("N80").Value = -7414.8
$ws.Range("H83").Value = 90913570
$ws.Range("I83").Value = 166670370
$ws.Range("J83").Value = 5418.8
$ws.Range("K83").Value = 833351850
$ws.Range("L83").Value = 27094
$ws.Range("M83").Value = -833346858
$ws.Range("N83").Value = -37078
$ws.Range("H133").Value = 90520.266
$ws.Range("J133").Value = 90200.28999999999
$ws.Range("L133").Value = 90200.28999999999
$ws.Range("N133").Value = -100320.29
$ws.Range("H138").Value = 106998
$ws.Range("J138").Value = 106998
$ws.Range("L138").Value = 106998
$ws.Range("N138").Value = -117278

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 462307.7
$ws.Range("I25").Value = 500000
$ws.Range("K25").Value = 500000
$ws.Range("M25").Value = -499770
$ws.Range("H46").Value = 3455.4167
$ws.Range("I46").Value = 2058.125
$ws.Range("K46").Value = 2058.125
$ws.Range("M46").Value = -1870.125
$ws.Range("H61").Value = 2307.9285
$ws.Range("I61").Value = 2242.5833
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 2242.5833
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -2040.5833
$ws.Range("N61").Value = -3104
$ws.Range("H100").Value = 9997.200000000001
$ws.Range("I100").Value = 9998.5
$ws.Range("J100").Value = 9996.333000000001
$ws.Range("K100").Value = 9998.5
$ws.Range("L100").Value = 9996.333000000001
$ws.Range("M100").Value = -9457.5
$ws.Range("N100").Value = -11078.333
$ws.Range("H113").Value = 2307.9285
$ws.Range("I113").Value = 2242.5833
$ws.Range("J113").Value = 2700
$ws.Range("K113").Value = 2242.5833
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -72.58329999999978
$ws.Range("N113").Value = -7040
$ws.Range("H136").Value = 4707.4546
$ws.Range("I136").Value = 1533.1765
$ws.Range("K136").Value = 4599.529500000001
$ws.Range("M136").Value = -2049.529500000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1221.6666
$ws.Range("I107").Value = 934.5238000000001
$ws.Range("K107").Value = 2803.5714
$ws.Range("M107").Value = -883.5714000000003
$ws.Range("H132").Value = 2425.7083
$ws.Range("I132").Value = 2180.0232
$ws.Range("K132").Value = 6540.069600000001
$ws.Range("M132").Value = -4010.069600000001

Write-Host "Applied market-data refresh: $($wb.Worksheets.Count) sheets processed"